$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.041.26'
$ws.Range("E2").Value = '  -2.96%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.865.46'
$ws.Range("E3").Value = '  -2.23%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.0000'
$ws.Range("E4").Value = '  -0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.92'
$ws.Range("E5").Value = '  -2.26%  '

# Row 6
$ws.Range("E6").Value = '  -0.05%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5160'
$ws.Range("E7").Value = '  -1.14%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3764'
$ws.Range("E8").Value = '  -0.74%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07161'
$ws.Range("E9").Value = '  -1.19%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8888'
$ws.Range("E10").Value = '  -2.44%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.66'
$ws.Range("E11").Value = '  -3.14%  '

# Row 12
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07608'
$ws.Range("E12").Value = '  -0.50%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.876.89'
$ws.Range("E13").Value = '  -2.24%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.300'
$ws.Range("E14").Value = '  -2.93%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.59'
$ws.Range("E15").Value = '  -2.97%  '

# Row 16
$ws.Range("E16").Value = '  -0.08%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008468'
$ws.Range("E17").Value = '  -2.83%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.07'
$ws.Range("E18").Value = '  -3.39%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("E19").Value = '  -0.03%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.069.48'
$ws.Range("E20").Value = '  -3.01%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.028'
$ws.Range("E21").Value = '  -2.50%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.132.92'

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.51'
$ws.Range("E23").Value = '  -3.24%  '

# Row 24
$ws.Range("E24").Value = '  -2.69%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.837'
$ws.Range("E25").Value = '  -1.75%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '147.59'
$ws.Range("E26").Value = '  -4.05%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.100'
$ws.Range("E28").Value = '  -3.40%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '112.76'
$ws.Range("E29").Value = '  -1.68%  '

# Row 30
$ws.Range("E30").Value = '  -4.42%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.681'
$ws.Range("E31").Value = '  -3.96%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09145'
$ws.Range("E32").Value = '  +1.43%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05124'
$ws.Range("E33").Value = '  -3.22%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.069'
$ws.Range("E34").Value = '  -3.57%  '

# Row 35
$ws.Range("E35").Value = '  -6.61%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7260'
$ws.Range("E36").Value = '  -7.19%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02035'
$ws.Range("E37").Value = '  -3.06%  '

# Row 38
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.062'
$ws.Range("E38").Value = '  -0.52%  '

# Row 39
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.497'
$ws.Range("E39").Value = '  -4.38%  '

# Row 40
$ws.Range("E40").Value = '  -1.98%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5330'
$ws.Range("E41").Value = '  -4.66%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.520'
$ws.Range("E42").Value = '  -3.23%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '116.43'
$ws.Range("E43").Value = '  +0.71%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.291'
$ws.Range("E44").Value = '  -3.34%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1466'
$ws.Range("E45").Value = '  -3.44%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4636'
$ws.Range("E46").Value = '  -3.93%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9994'
$ws.Range("E47").Value = '  -0.06%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.981'
$ws.Range("E48").Value = '  -4.71%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.571'
$ws.Range("E49").Value = '  -3.22%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.53'
$ws.Range("E50").Value = '  -1.47%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.59'
$ws.Range("E51").Value = '  -5.27%  '
